# Apply updated result values (row 2) across the yearly result sheets.
# Values were recomputed "from server" per the commit message.

$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = @{
        "N2" = 7098.474579830174
        "O2" = 6965.1427089074
    }
    "2030" = @{
        "B2" = 5707.815717280662
        "I2" = 44492.05901988943
        "L2" = 66334.06707325629
        "M2" = 21991.42050229464
        "N2" = 10503.58911838768
        "O2" = 12039.81600698566
    }
    "2035" = @{
        "A2" = 2927.360317916481
        "B2" = 7940.887964949257
        "E2" = 67179.99183625776
        "I2" = 59530.75343380851
        "L2" = 66334.06707325629
        "M2" = 25547.11936466757
        "N2" = 14996.33602230955
        "O2" = 14718.4361623581
    }
    "2040" = @{
        "A2" = 2927.360317916481
        "B2" = 7940.887964949257
        "E2" = 67179.99183625776
        "I2" = 59530.75343380851
        "L2" = 66334.06707325629
        "M2" = 25547.11936466757
        "N2" = 15101.21309504717
        "O2" = 14718.4361623581
    }
    "2045" = @{
        "A2" = 6352.985609279765
        "B2" = 7940.887964949257
        "E2" = 67179.99183625776
        "I2" = 59530.75343380851
        "L2" = 66334.06707325629
        "M2" = 25547.11936466757
        "N2" = 15632.92561905594
        "O2" = 17025.5453334183
    }
    "2050" = @{
        "A2" = 6352.985609279765
        "B2" = 7940.887964949257
        "E2" = 67179.99183625776
        "I2" = 59530.75343380851
        "L2" = 66334.06707325629
        "M2" = 25547.11936466757
        "N2" = 15632.92561905594
        "O2" = 17025.5453334183
    }
}

foreach ($sheetNameKey in $updates.Keys) {
    $sheetName = [string]$sheetNameKey
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetNameKey]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
